$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.085101972571708
$ws.Cells.Item(2, 4).Value = 1.076127538284687
$ws.Cells.Item(2, 5).Value = 1.09771092208919
$ws.Cells.Item(2, 6).Value = 1.102347599429943
$ws.Cells.Item(2, 9).Value = 1.05027148432415
$ws.Cells.Item(2, 10).Value = 1.089958490888531
$ws.Cells.Item(2, 11).Value = 1.078812289038594
$ws.Cells.Item(2, 12).Value = 1.100339972706455
$ws.Cells.Item(2, 13).Value = 1.104964990862307
$ws.Cells.Item(2, 14).Value = 1.091506356965579
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.087320635992203
$ws.Cells.Item(3, 4).Value = 1.077847353611072
$ws.Cells.Item(3, 5).Value = 1.099899771687645
$ws.Cells.Item(3, 6).Value = 1.10454703164576
$ws.Cells.Item(3, 9).Value = 1.050882667152107
$ws.Cells.Item(3, 10).Value = 1.091834451802705
$ws.Cells.Item(3, 11).Value = 1.080347343488697
$ws.Cells.Item(3, 12).Value = 1.102347057972067
$ws.Cells.Item(3, 13).Value = 1.10698349928836
$ws.Cells.Item(3, 14).Value = 1.093384981959427
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.088750260508363
$ws.Cells.Item(4, 4).Value = 1.078954610733508
$ws.Cells.Item(4, 5).Value = 1.101310432426495
$ws.Cells.Item(4, 6).Value = 1.105964536251188
$ws.Cells.Item(4, 9).Value = 1.051273979814478
$ws.Cells.Item(4, 10).Value = 1.093042158161293
$ws.Cells.Item(4, 11).Value = 1.081334590251575
$ws.Cells.Item(4, 12).Value = 1.103639729384229
$ws.Cells.Item(4, 13).Value = 1.108283557506727
$ws.Cells.Item(4, 14).Value = 1.09459440339957
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.089349872293312
$ws.Cells.Item(5, 4).Value = 1.079418791778182
$ws.Cells.Item(5, 5).Value = 1.101902148863568
$ws.Cells.Item(5, 6).Value = 1.106559129045453
$ws.Cells.Item(5, 9).Value = 1.051437501754988
$ws.Cells.Item(5, 10).Value = 1.093548430943893
$ws.Cells.Item(5, 11).Value = 1.081748208062571
$ws.Cells.Item(5, 12).Value = 1.1041817513921
$ws.Cells.Item(5, 13).Value = 1.108828683911388
$ws.Cells.Item(5, 14).Value = 1.095101395147592
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.089450468492413
$ws.Cells.Item(6, 4).Value = 1.079496653748069
$ws.Cells.Item(6, 5).Value = 1.10200142387411
$ws.Cells.Item(6, 6).Value = 1.106658886964871
$ws.Cells.Item(6, 9).Value = 1.051464900307689
$ws.Cells.Item(6, 10).Value = 1.093633352362178
$ws.Cells.Item(6, 11).Value = 1.081817573703467
$ws.Cells.Item(6, 12).Value = 1.104272677059401
$ws.Cells.Item(6, 13).Value = 1.108920130764398
$ws.Cells.Item(6, 14).Value = 1.095186437164031
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.088758278015277
$ws.Cells.Item(7, 4).Value = 1.078960818255682
$ws.Cells.Item(7, 5).Value = 1.101318344134644
$ws.Cells.Item(7, 6).Value = 1.105972486396383
$ws.Cells.Item(7, 9).Value = 1.05127616866195
$ws.Cells.Item(7, 10).Value = 1.093048928649446
$ws.Cells.Item(7, 11).Value = 1.081340122582778
$ws.Cells.Item(7, 12).Value = 1.103646977435758
$ws.Cells.Item(7, 13).Value = 1.108290847043193
$ws.Cells.Item(7, 14).Value = 1.094601183502593
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.085853042909693
$ws.Cells.Item(8, 4).Value = 1.076709928502908
$ws.Cells.Item(8, 5).Value = 1.098451847967939
$ws.Cells.Item(8, 6).Value = 1.103092102602441
$ws.Cells.Item(8, 9).Value = 1.050478905077554
$ws.Cells.Item(8, 10).Value = 1.090593774986839
$ws.Cells.Item(8, 11).Value = 1.079332332440152
$ws.Cells.Item(8, 12).Value = 1.101019547177929
$ws.Cells.Item(8, 13).Value = 1.105648426898199
$ws.Cells.Item(8, 14).Value = 1.092142543240175
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.080686168704518
$ws.Cells.Item(9, 4).Value = 1.07269968736428
$ws.Cells.Item(9, 5).Value = 1.093355812279202
$ws.Cells.Item(9, 6).Value = 1.097971561376338
$ws.Cells.Item(9, 9).Value = 1.049041645189526
$ws.Cells.Item(9, 10).Value = 1.086218937338253
$ws.Cells.Item(9, 11).Value = 1.07574701709059
$ws.Cells.Item(9, 12).Value = 1.096341979827865
$ws.Cells.Item(9, 13).Value = 1.100944403760877
$ws.Cells.Item(9, 14).Value = 1.087761492820327
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.077207501509894
$ws.Cells.Item(10, 4).Value = 1.069995046821423
$ws.Cells.Item(10, 5).Value = 1.089926194372177
$ws.Cells.Item(10, 6).Value = 1.094525577128141
$ws.Cells.Item(10, 9).Value = 1.04806100616259
$ws.Cells.Item(10, 10).Value = 1.083267877758757
$ws.Cells.Item(10, 11).Value = 1.073323425778192
$ws.Cells.Item(10, 12).Value = 1.093189574187706
$ws.Cells.Item(10, 13).Value = 1.0977743253434
$ws.Cells.Item(10, 14).Value = 1.084806242397737
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.075692629223227
$ws.Cells.Item(11, 4).Value = 1.068816150654583
$ws.Cells.Item(11, 5).Value = 1.088433014146284
$ws.Cells.Item(11, 6).Value = 1.093025301658855
$ws.Cells.Item(11, 9).Value = 1.04763089312032
$ws.Cells.Item(11, 10).Value = 1.081981433604305
$ws.Cells.Item(11, 11).Value = 1.072265715598458
$ws.Cells.Item(11, 12).Value = 1.091816036350159
$ws.Cells.Item(11, 13).Value = 1.096393124862868
$ws.Cells.Item(11, 14).Value = 1.083517971345028
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.075128605933752
$ws.Cells.Item(12, 4).Value = 1.068377056807712
$ws.Cells.Item(12, 5).Value = 1.087877118159009
$ws.Cells.Item(12, 6).Value = 1.09246676877537
$ws.Cells.Item(12, 9).Value = 1.047470291506783
$ws.Cells.Item(12, 10).Value = 1.081502260064598
$ws.Cells.Item(12, 11).Value = 1.07187156057857
$ws.Cells.Item(12, 12).Value = 1.091304524120387
$ws.Cells.Item(12, 13).Value = 1.095878764705094
$ws.Cells.Item(12, 14).Value = 1.083038117323933
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.075249651818349
$ws.Cells.Item(13, 4).Value = 1.068471298721851
$ws.Cells.Item(13, 5).Value = 1.087996417517319
$ws.Cells.Item(13, 6).Value = 1.092586633822617
$ws.Cells.Item(13, 9).Value = 1.047504779263336
$ws.Cells.Item(13, 10).Value = 1.081605105233737
$ws.Cells.Item(13, 11).Value = 1.071956166331634
$ws.Cells.Item(13, 12).Value = 1.091414305544611
$ws.Cells.Item(13, 13).Value = 1.095989157098001
$ws.Cells.Item(13, 14).Value = 1.083141108545008
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.075646034307278
$ws.Cells.Item(14, 4).Value = 1.068779879673875
$ws.Cells.Item(14, 5).Value = 1.088387089599172
$ws.Cells.Item(14, 6).Value = 1.092979159174529
$ws.Cells.Item(14, 9).Value = 1.047617634925636
$ws.Cells.Item(14, 10).Value = 1.0819418523333
$ws.Cells.Item(14, 11).Value = 1.072233160796711
$ws.Cells.Item(14, 12).Value = 1.091773781723202
$ws.Cells.Item(14, 13).Value = 1.096350634860446
$ws.Cells.Item(14, 14).Value = 1.083478333864078
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.075890080824916
$ws.Cells.Item(15, 4).Value = 1.068969846690537
$ws.Cells.Item(15, 5).Value = 1.088627627063665
$ws.Cells.Item(15, 6).Value = 1.093220838309571
$ws.Cells.Item(15, 9).Value = 1.047687057503558
$ws.Cells.Item(15, 10).Value = 1.082149155876027
$ws.Cells.Item(15, 11).Value = 1.072403656440634
$ws.Cells.Item(15, 12).Value = 1.091995090908085
$ws.Cells.Item(15, 13).Value = 1.096573177059999
$ws.Cells.Item(15, 14).Value = 1.083685931801614
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.077307854098073
$ws.Cells.Item(16, 4).Value = 1.070073119865491
$ws.Cells.Item(16, 5).Value = 1.090025117034085
$ws.Cells.Item(16, 6).Value = 1.09462497048778
$ws.Cells.Item(16, 9).Value = 1.048089434520652
$ws.Cells.Item(16, 10).Value = 1.08335307019607
$ws.Cells.Item(16, 11).Value = 1.073393445494802
$ws.Cells.Item(16, 12).Value = 1.093280548498355
$ws.Cells.Item(16, 13).Value = 1.097865807970763
$ws.Cells.Item(16, 14).Value = 1.084891555818082
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.078194857383044
$ws.Cells.Item(17, 4).Value = 1.070763071509787
$ws.Cells.Item(17, 5).Value = 1.090899519777784
$ws.Cells.Item(17, 6).Value = 1.095503537344602
$ws.Cells.Item(17, 9).Value = 1.048340355612393
$ws.Cells.Item(17, 10).Value = 1.084105921586004
$ws.Cells.Item(17, 11).Value = 1.074012076082223
$ws.Cells.Item(17, 12).Value = 1.094084572783395
$ws.Cells.Item(17, 13).Value = 1.098674329052839
$ws.Cells.Item(17, 14).Value = 1.085645476343344
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.078711405816565
$ws.Cells.Item(18, 4).Value = 1.071164761236401
$ws.Cells.Item(18, 5).Value = 1.091408762152629
$ws.Cells.Item(18, 6).Value = 1.096015207763886
$ws.Cells.Item(18, 9).Value = 1.048486184883825
$ws.Cells.Item(18, 10).Value = 1.084544218251522
$ws.Cells.Item(18, 11).Value = 1.074372116107896
$ws.Cells.Item(18, 12).Value = 1.094552726587759
$ws.Cells.Item(18, 13).Value = 1.099145104787371
$ws.Cells.Item(18, 14).Value = 1.0860843954404
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.07888739671956
$ws.Cells.Item(19, 4).Value = 1.071301601245152
$ws.Cells.Item(19, 5).Value = 1.09158226922414
$ws.Cells.Item(19, 6).Value = 1.096189542601885
$ws.Cells.Item(19, 9).Value = 1.048535819621116
$ws.Cells.Item(19, 10).Value = 1.084693526600763
$ws.Cells.Item(19, 11).Value = 1.074494746228679
$ws.Cells.Item(19, 12).Value = 1.09471221700616
$ws.Cells.Item(19, 13).Value = 1.099305489055155
$ws.Cells.Item(19, 14).Value = 1.086233915824623
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.078099776021926
$ws.Cells.Item(20, 4).Value = 1.070689123765358
$ws.Cells.Item(20, 5).Value = 1.090805785810445
$ws.Cells.Item(20, 6).Value = 1.095409356692277
$ws.Cells.Item(20, 9).Value = 1.048313488944043
$ws.Cells.Item(20, 10).Value = 1.084025233746835
$ws.Cells.Item(20, 11).Value = 1.073945785503638
$ws.Cells.Item(20, 12).Value = 1.093998393674007
$ws.Cells.Item(20, 13).Value = 1.098587667581639
$ws.Cells.Item(20, 14).Value = 1.085564673918187
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.075529346664147
$ws.Cells.Item(21, 4).Value = 1.068689043642566
$ws.Cells.Item(21, 5).Value = 1.088272081584365
$ws.Cells.Item(21, 6).Value = 1.09286360545776
$ws.Cells.Item(21, 9).Value = 1.047584425003205
$ws.Cells.Item(21, 10).Value = 1.081842725741099
$ws.Cells.Item(21, 11).Value = 1.072151628257106
$ws.Cells.Item(21, 12).Value = 1.091667961673217
$ws.Cells.Item(21, 13).Value = 1.096244225441504
$ws.Cells.Item(21, 14).Value = 1.083379066500746
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.073905486113721
$ws.Cells.Item(22, 4).Value = 1.0674245586628
$ws.Cells.Item(22, 5).Value = 1.086671716939462
$ws.Cells.Item(22, 6).Value = 1.091255658175646
$ws.Cells.Item(22, 9).Value = 1.047121175426844
$ws.Cells.Item(22, 10).Value = 1.080462777614424
$ws.Cells.Item(22, 11).Value = 1.071016181793307
$ws.Cells.Item(22, 12).Value = 1.090195075104456
$ws.Cells.Item(22, 13).Value = 1.094763149080453
$ws.Cells.Item(22, 14).Value = 1.081997158689449
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.074767072295433
$ws.Cells.Item(23, 4).Value = 1.068095556223669
$ws.Cells.Item(23, 5).Value = 1.087520808468652
$ws.Cells.Item(23, 6).Value = 1.092108770225386
$ws.Cells.Item(23, 9).Value = 1.047367217922213
$ws.Cells.Item(23, 10).Value = 1.081195058213999
$ws.Cells.Item(23, 11).Value = 1.071618814130929
$ws.Cells.Item(23, 12).Value = 1.090976618370037
$ws.Cells.Item(23, 13).Value = 1.095549034888502
$ws.Cells.Item(23, 14).Value = 1.082730479211469
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.078142741745585
$ws.Cells.Item(24, 4).Value = 1.070722539876213
$ws.Cells.Item(24, 5).Value = 1.090848142570163
$ws.Cells.Item(24, 6).Value = 1.095451915292219
$ws.Cells.Item(24, 9).Value = 1.048325630468324
$ws.Cells.Item(24, 10).Value = 1.084061695667741
$ws.Cells.Item(24, 11).Value = 1.073975741823325
$ws.Cells.Item(24, 12).Value = 1.094037336836128
$ws.Cells.Item(24, 13).Value = 1.098626828705309
$ws.Cells.Item(24, 14).Value = 1.085601187619204
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.082027776432071
$ws.Cells.Item(25, 4).Value = 1.073741798674359
$ws.Cells.Item(25, 5).Value = 1.094678791962751
$ws.Cells.Item(25, 6).Value = 1.099300881005407
$ws.Cells.Item(25, 9).Value = 1.049417119337198
$ws.Cells.Item(25, 10).Value = 1.087355880284098
$ws.Cells.Item(25, 11).Value = 1.07667967491328
$ws.Cells.Item(25, 12).Value = 1.097557095216537
$ws.Cells.Item(25, 13).Value = 1.102166364316597
$ws.Cells.Item(25, 14).Value = 1.088900050355564
